$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "emailVerified" column (column E) is no longer needed since the
# email field itself is now optional. Remove the entire column, which
# shifts phone/gender/studentId (F,G,H) left into E,F,G.
$ws.Columns.Item(5).Delete()

# The "email" column description (now column C, row 2) changes from
# "required, email must not be empty" to "optional, email".
$ws.Cells.Item(2, 3).Value = "选填，邮箱"

# Update the active selection to match the target state.
$ws.Range("D2").Select()
